$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").EntireColumn.Insert()
$ws.Range("D1").EntireColumn.Insert()
$ws.Range("F1").EntireColumn.Insert()

$ws.Range("B1").Value = "ProductQty1"
$ws.Range("D1").Value = "ProductQty2"
$ws.Range("F1").Value = "ProductQty3"

$ws.Range("B2").Value = 1
$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("F2").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("F5").Value = 1
$ws.Range("F6").Value = 1
$ws.Range("F7").Value = 1

$ws.Columns.Item(2).ColumnWidth = 14
$ws.Columns.Item(4).ColumnWidth = 14
$ws.Columns.Item(6).ColumnWidth = 14

$ws.Range("F3:F7").Select() | Out-Null

